$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix hardcoded test data values in column F (Test Case Title)
$ws.Range("F2").Value = "loginTest"
$ws.Range("F3").Value = "loginTest"
$ws.Range("F4").Value = "Homepage"

# Update selection / view to F3 (and reset scrolled top-left cell to default)
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("F3").Select()
